# Word COM-interop edit script
# Commit message: "updated date on updated cv"
#
# The CV's "last updated" footer date moves from "August 4, 2018" to
# "October 3, 2018". While touching that part of the document the CV's
# bibliography entry for "Working with NHANES data in R..." also ends
# up re-flowed into a single run (same visible text, just no longer
# split across two runs).

$d = $word.ActiveDocument

# --- 1) Update the footer date -------------------------------------------
$footer = $d.Sections(1).Footers(1)
$dateFound = $footer.Range.Find.Execute(
    "August 4, 2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "October 3, 2018", 2)
if (-not $dateFound) {
    Write-Output "WARNING: footer date 'August 4, 2018' not found"
}

# --- 2) Re-flow the NHANES bullet so its text lives in one run -----------
$nhanesFound = $d.Content.Find.Execute(
    "Information Systems, Systems Biology", $true, $false, $false, $false,
    $false, $true, 1, $false, "Information Systems, Systems Biology", 2)
if (-not $nhanesFound) {
    Write-Output "WARNING: NHANES bullet text not found"
}
